# Update the handback/handoff timestamps for the first data row ("de1ff3f0-...")
# on both the "zh-cn" and "de-de" worksheets, as part of regenerating the
# handback status report.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-13 15:06:31"
$wsZhCn.Range("H2").Value = "2016-03-13 15:06:48"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-13 15:06:35"
$wsDeDe.Range("H2").Value = "2016-03-13 15:06:54"
